# Nifty 50 Return of Change Tracker - update summary sheet to 30-Day ROC
# Updates the header label, recomputes the 30-Day ROC (%) values, and
# re-sorts all data rows (2-51) in descending order of the new ROC value.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header for column C
$ws.Range("C1").Value = "30-Day ROC (%)"

# Ticker -> new 30-Day ROC (%) value (Latest Close in column B is unchanged per ticker)
$rocMap = [ordered]@{
    "ETERNAL.NS"    = 22.06248282399144
    "HEROMOTOCO.NS" = 9.221453881006681
    "HINDUNILVR.NS" = 7.227466169731667
    "ASIANPAINT.NS" = 4.02320345374283
    "CIPLA.NS"      = 3.66541065936743
    "APOLLOHOSP.NS" = 3.390614672835435
    "M&M.NS"        = 2.866683555931693
    "LT.NS"         = 2.62636691565663
    "SBIN.NS"       = 2.409864055925781
    "SBILIFE.NS"    = 1.679373839027143
    "MARUTI.NS"     = 1.442910915934759
    "NTPC.NS"       = 1.419605473208052
    "EICHERMOT.NS"  = 0.8397480755773179
    "JIOFIN.NS"     = 0.7229676024514742
    "HDFCLIFE.NS"   = 0.3690302862070061
    "HINDALCO.NS"   = 0.3463988595846956
    "HDFCBANK.NS"   = 0.2012983769888343
    "ICICIBANK.NS"  = 0.09818531384329088
    "JSWSTEEL.NS"   = -0.05735356083929055
    "ITC.NS"        = -0.459657401337632
    "COALINDIA.NS"  = -0.5175313589920516
    "ULTRACEMCO.NS" = -0.6533311824487797
    "POWERGRID.NS"  = -1.702417361997066
    "GRASIM.NS"     = -1.870393011284222
    "BAJAJ-AUTO.NS" = -2.016104980614375
    "SUNPHARMA.NS"  = -2.28016478366222
    "DRREDDY.NS"    = -2.392939150130502
    "BAJAJFINSV.NS" = -2.807087239692307
    "ONGC.NS"       = -2.913337644516922
    "TATACONSUM.NS" = -3.506523993583854
    "TATAMOTORS.NS" = -3.736971019087276
    "TITAN.NS"      = -5.145840958510961
    "BAJFINANCE.NS" = -5.350767554874636
    "TATASTEEL.NS"  = -6.395033678595885
    "KOTAKBANK.NS"  = -6.978281748784287
    "BHARTIARTL.NS" = -7.205467066889792
    "WIPRO.NS"      = -7.579101726103099
    "SHRIRAMFIN.NS" = -8.728162581600852
    "NESTLEIND.NS"  = -8.775743681363558
    "AXISBANK.NS"   = -8.786613146548694
    "ADANIPORTS.NS" = -9.139817552177576
    "RELIANCE.NS"   = -9.487415691624246
    "BEL.NS"        = -9.700881197305723
    "INFY.NS"       = -10.56403318454434
    "INDUSINDBK.NS" = -10.74327710936829
    "TECHM.NS"      = -11.07722175595357
    "TCS.NS"        = -11.12973402039401
    "ADANIENT.NS"   = -12.62254364579414
    "HCLTECH.NS"    = -12.97388705055696
    "TRENT.NS"      = -13.25310935228558
}

# Read the existing data rows (2-51): Ticker, Latest Close
# NOTE: use .Value2 for reads - .Value is unreliable for round-tripping
# already-populated cells in this COM host.
$firstRow = 2
$lastRow = 51
$rows = @()
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ticker = $ws.Cells.Item($r, 1).Value2
    $close = $ws.Cells.Item($r, 2).Value2
    $roc = $rocMap[$ticker]
    $rows += @{ Ticker = $ticker; Close = $close; Roc = $roc }
}

# Sort rows descending by the new ROC value
$sorted = $rows | Sort-Object -Property Roc -Descending

# Write the sorted rows back
$r = $firstRow
foreach ($row in $sorted) {
    $ws.Cells.Item($r, 1).Value = $row.Ticker
    $ws.Cells.Item($r, 2).Value = $row.Close
    $ws.Cells.Item($r, 3).Value = $row.Roc
    $r++
}
